# Add a new "2022-Q3" quarterly sheet (with fund-holding detail data) to the
# workbook, positioned right after "总计" and before "2022-Q2", and add a
# matching summary row to the "总计" (total) sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Build the new "2022-Q3" worksheet by copying the existing "2022-Q2"
#    sheet (same column layout/styles) immediately before it, then
#    overwriting its data with the 2022-Q3 numbers.
# ---------------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2022-Q2")
$templateSheet.Copy($templateSheet)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Fund holding detail rows for 2022-Q3:
# index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$rows = @(
    @(0,  "506007", "广发科创板两年定开混合",         "5.01", "94.25", "4.98", "0.2495", 6),
    @(1,  "009414", "中银大健康股票A",                "2.64", "87.62", "3.38", "0.0892", 7),
    @(2,  "002801", "泓德泓信灵活配置混合",           "2.20", "92.10", "3.08", "0.0678", 7),
    @(3,  "001695", "泓德泓业灵活配置混合",           "1.21", "89.73", "4.50", "0.0544", 7),
    @(4,  "011781", "泓德慧享混合A",                  "3.73", "26.63", "0.82", "0.0306", 7),
    @(5,  "009015", "泓德睿享一年持有期混合A",        "2.60", "27.82", "1.16", "0.0302", 6),
    @(6,  "002681", "金鹰元和灵活配置混合A",          "0.30", "87.46", "4.64", "0.0139", 9),
    @(7,  "002682", "金鹰元和灵活配置混合C",          "0.23", "87.46", "4.64", "0.0107", 9),
    @(8,  "010321", "中银大健康股票C",                "0.11", "87.62", "3.38", "0.0037", 7),
    @(9,  "009016", "泓德睿享一年持有期混合C",        "0.06", "27.82", "1.16", "0.0007", 6),
    @(10, "011782", "泓德慧享混合C",                  "0.00", "26.63", "0.82", "0",      7)
)
$lastRow = 1 + $rows.Count   # 12

# Columns B (fund code) and D:G (numeric-looking figures) are stored as
# *text* in the source data. Force text storage for those blocks up front so
# assigning "506007" / "5.01" etc. doesn't get auto-coerced to a number
# (which would also strip the leading zeros off codes like "009414").
$q3.Range("B2:B" + $lastRow).NumberFormat = "@"
$q3.Range("D2:G" + $lastRow).NumberFormat = "@"

# The very last row's G value (持有市值) is a literal number 0 in the source,
# not text "0" - leave it in the default (General) style.
$q3.Range("G" + $lastRow).ClearFormats()

# Rows 2-9 already exist (copied from the 2022-Q2 template) with column A
# bold/bordered (style matches the header). Rows 10-12 are brand new, so
# column A needs that same look copied over explicitly.
$q3.Range("A9").Copy()
$q3.Range("A10:A" + $lastRow).PasteSpecial(-4122)
$q3.Application.CutCopyMode = $false

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# G12 (last row's 持有市值) is a real number 0, not text "0"
$q3.Range("G" + $lastRow).Value = 0

$wb.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Insert the 2022-Q3 summary row into the "总计" sheet (row 2), pushing
#    the existing quarters down, and renumber the index column.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# Match the bold/bordered style used by the other index cells in column A
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Application.CutCopyMode = $false

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 11
$total.Range("D2").Value = 0.55

# Renumber the remaining rows' index column sequentially (1..7)
for ($row = 3; $row -le 9; $row++) {
    $total.Cells.Item($row, 1).Value = $row - 2
}
